# Auto-generated: apply market-price/profit refresh to Mateus_Profits workbook
# Updates columns H:N (currentAveragePrice.. / LevePrice.. / LeveProfit..) for the
# rows whose underlying market data changed, per sheet (ALC, ARM, BSM, CRP, CUL, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 35714440
$ws.Cells.Item(33, 9).Value = 41666824
$ws.Cells.Item(33, 10).Value = 145
$ws.Cells.Item(33, 11).Value = 41666824
$ws.Cells.Item(33, 12).Value = 145
$ws.Cells.Item(33, 13).Value = -41666595
$ws.Cells.Item(33, 14).Value = -603
$ws.Cells.Item(51, 8).Value = 13827.333
$ws.Cells.Item(51, 9).Value = 13624.5
$ws.Cells.Item(51, 10).Value = 13885.286
$ws.Cells.Item(51, 11).Value = 13624.5
$ws.Cells.Item(51, 12).Value = 13885.286
$ws.Cells.Item(51, 13).Value = -13140.5
$ws.Cells.Item(51, 14).Value = -14853.286
$ws.Cells.Item(58, 8).Value = 288
$ws.Cells.Item(58, 9).Value = 288
$ws.Cells.Item(58, 11).Value = 864
$ws.Cells.Item(58, 13).Value = -714
$ws.Cells.Item(86, 8).Value = 2169.318
$ws.Cells.Item(86, 9).Value = 2340.3333
$ws.Cells.Item(86, 11).Value = 2340.3333
$ws.Cells.Item(86, 13).Value = -1217.3333
$ws.Cells.Item(89, 8).Value = 2169.318
$ws.Cells.Item(89, 9).Value = 2340.3333
$ws.Cells.Item(89, 11).Value = 11701.6665
$ws.Cells.Item(89, 13).Value = -6085.666499999999
$ws.Cells.Item(113, 8).Value = 4161.263
$ws.Cells.Item(113, 9).Value = 3318.8
$ws.Cells.Item(113, 10).Value = 5097.3335
$ws.Cells.Item(113, 11).Value = 3318.8
$ws.Cells.Item(113, 12).Value = 5097.3335
$ws.Cells.Item(113, 13).Value = -64.80000000000018
$ws.Cells.Item(113, 14).Value = -11605.3335
$ws.Cells.Item(116, 8).Value = 3586
$ws.Cells.Item(116, 9).Value = 2496.3333
$ws.Cells.Item(116, 10).Value = 4053
$ws.Cells.Item(116, 11).Value = 2496.3333
$ws.Cells.Item(116, 12).Value = 4053
$ws.Cells.Item(116, 13).Value = 945.6667000000002
$ws.Cells.Item(116, 14).Value = -10937
$ws.Cells.Item(132, 8).Value = 7673.375
$ws.Cells.Item(132, 9).Value = 6609.3687
$ws.Cells.Item(132, 10).Value = 9228.462
$ws.Cells.Item(132, 11).Value = 19828.1061
$ws.Cells.Item(132, 12).Value = 27685.386
$ws.Cells.Item(132, 13).Value = -17298.1061
$ws.Cells.Item(132, 14).Value = -32745.386
$ws.Cells.Item(138, 8).Value = 2319.7407
$ws.Cells.Item(138, 9).Value = 1507.9286
$ws.Cells.Item(138, 10).Value = 3194
$ws.Cells.Item(138, 11).Value = 4523.7858
$ws.Cells.Item(138, 12).Value = 9582
$ws.Cells.Item(138, 13).Value = 616.2142000000003
$ws.Cells.Item(138, 14).Value = -19862

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1960.7764
$ws.Cells.Item(32, 9).Value = 1960.7764
$ws.Cells.Item(32, 11).Value = 1960.7764
$ws.Cells.Item(32, 13).Value = -1673.7764
$ws.Cells.Item(74, 8).Value = 5365.24
$ws.Cells.Item(74, 9).Value = 3796.3684
$ws.Cells.Item(74, 11).Value = 3796.3684
$ws.Cells.Item(74, 13).Value = -2922.3684
$ws.Cells.Item(77, 8).Value = 5365.24
$ws.Cells.Item(77, 9).Value = 3796.3684
$ws.Cells.Item(77, 11).Value = 18981.842
$ws.Cells.Item(77, 13).Value = -14613.842
$ws.Cells.Item(97, 8).Value = 569.17645
$ws.Cells.Item(97, 9).Value = 587.6129
$ws.Cells.Item(97, 11).Value = 587.6129
$ws.Cells.Item(97, 13).Value = -91.61289999999997
$ws.Cells.Item(122, 8).Value = 1410.0625
$ws.Cells.Item(122, 9).Value = 1183.2142
$ws.Cells.Item(122, 11).Value = 3549.6426
$ws.Cells.Item(122, 13).Value = -1099.6426
$ws.Cells.Item(132, 8).Value = 2526.795
$ws.Cells.Item(132, 9).Value = 2461.7104
$ws.Cells.Item(132, 11).Value = 7385.1312
$ws.Cells.Item(132, 13).Value = -4855.1312
$ws.Cells.Item(135, 8).Value = 172500
$ws.Cells.Item(135, 10).Value = 172500
$ws.Cells.Item(135, 12).Value = 172500
$ws.Cells.Item(135, 14).Value = -182640

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 6750
$ws.Cells.Item(81, 10).Value = 6750
$ws.Cells.Item(81, 12).Value = 6750
$ws.Cells.Item(81, 14).Value = -8872
$ws.Cells.Item(84, 8).Value = 6750
$ws.Cells.Item(84, 10).Value = 6750
$ws.Cells.Item(84, 12).Value = 20250
$ws.Cells.Item(84, 14).Value = -30858
$ws.Cells.Item(97, 8).Value = 18325
$ws.Cells.Item(97, 9).Value = 16984
$ws.Cells.Item(97, 10).Value = 21007
$ws.Cells.Item(97, 11).Value = 16984
$ws.Cells.Item(97, 12).Value = 21007
$ws.Cells.Item(97, 13).Value = -15993
$ws.Cells.Item(97, 14).Value = -22989
$ws.Cells.Item(134, 8).Value = 2759
$ws.Cells.Item(134, 9).Value = 2745.4443
$ws.Cells.Item(134, 10).Value = 3125
$ws.Cells.Item(134, 11).Value = 8236.332900000001
$ws.Cells.Item(134, 12).Value = 9375
$ws.Cells.Item(134, 13).Value = -5701.332900000001
$ws.Cells.Item(134, 14).Value = -14445
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 10000
$ws.Cells.Item(17, 9).Value = 10000
$ws.Cells.Item(17, 11).Value = 10000
$ws.Cells.Item(17, 13).Value = -9826
$ws.Cells.Item(25, 8).Value = 11055.5
$ws.Cells.Item(25, 10).Value = 11000
$ws.Cells.Item(25, 12).Value = 11000
$ws.Cells.Item(25, 14).Value = -11348
$ws.Cells.Item(31, 8).Value = 7826.154
$ws.Cells.Item(31, 9).Value = 6999.6665
$ws.Cells.Item(31, 11).Value = 6999.6665
$ws.Cells.Item(31, 13).Value = -6704.6665
$ws.Cells.Item(34, 8).Value = 7826.154
$ws.Cells.Item(34, 9).Value = 6999.6665
$ws.Cells.Item(34, 11).Value = 6999.6665
$ws.Cells.Item(34, 13).Value = -6797.6665
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).Value = $null
$ws.Cells.Item(58, 8).Value = 6521.2593
$ws.Cells.Item(58, 9).Value = 4691.6
$ws.Cells.Item(58, 10).Value = 8808.333000000001
$ws.Cells.Item(58, 11).Value = 4691.6
$ws.Cells.Item(58, 12).Value = 8808.333000000001
$ws.Cells.Item(58, 13).Value = -4488.6
$ws.Cells.Item(58, 14).Value = -9214.333000000001
$ws.Cells.Item(62, 8).Value = 6999.6
$ws.Cells.Item(62, 9).Value = 6500
$ws.Cells.Item(62, 10).Value = 7124.5
$ws.Cells.Item(62, 11).Value = 6500
$ws.Cells.Item(62, 12).Value = 7124.5
$ws.Cells.Item(62, 13).Value = -5876
$ws.Cells.Item(62, 14).Value = -8372.5
$ws.Cells.Item(65, 8).Value = 6999.6
$ws.Cells.Item(65, 9).Value = 6500
$ws.Cells.Item(65, 10).Value = 7124.5
$ws.Cells.Item(65, 11).Value = 32500
$ws.Cells.Item(65, 12).Value = 35622.5
$ws.Cells.Item(65, 13).Value = -29380
$ws.Cells.Item(65, 14).Value = -41862.5
$ws.Cells.Item(132, 8).Value = 4641.909
$ws.Cells.Item(132, 9).Value = 4917.8887
$ws.Cells.Item(132, 11).Value = 14753.6661
$ws.Cells.Item(132, 13).Value = -12223.6661
$ws.Cells.Item(136, 8).Value = 6521.2593
$ws.Cells.Item(136, 9).Value = 4691.6
$ws.Cells.Item(136, 10).Value = 8808.333000000001
$ws.Cells.Item(136, 11).Value = 14074.8
$ws.Cells.Item(136, 12).Value = 26424.999
$ws.Cells.Item(136, 13).Value = -11524.8
$ws.Cells.Item(136, 14).Value = -31524.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 137423.64
$ws.Cells.Item(37, 10).Value = 137423.64
$ws.Cells.Item(37, 12).Value = 412270.92
$ws.Cells.Item(37, 14).Value = -412494.92
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).Value = $null
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).Value = $null
$ws.Cells.Item(128, 8).Value = 2174091
$ws.Cells.Item(128, 9).Value = 2174091
$ws.Cells.Item(128, 11).Value = 6522273
$ws.Cells.Item(128, 13).Value = -6517293

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1497.4166
$ws.Cells.Item(22, 10).Value = 1468
$ws.Cells.Item(22, 12).Value = 1468
$ws.Cells.Item(22, 14).Value = -2058
$ws.Cells.Item(27, 8).Value = 1497.4166
$ws.Cells.Item(27, 10).Value = 1468
$ws.Cells.Item(27, 12).Value = 1468
$ws.Cells.Item(27, 14).Value = -1682
$ws.Cells.Item(61, 8).Value = 102196.3
$ws.Cells.Item(61, 10).Value = 1499.5
$ws.Cells.Item(61, 12).Value = 1499.5
$ws.Cells.Item(61, 14).Value = -1903.5
$ws.Cells.Item(113, 8).Value = 102196.3
$ws.Cells.Item(113, 10).Value = 1499.5
$ws.Cells.Item(113, 12).Value = 1499.5
$ws.Cells.Item(113, 14).Value = -5839.5
$ws.Cells.Item(136, 8).Value = 3827.925
$ws.Cells.Item(136, 9).Value = 3714.8857
$ws.Cells.Item(136, 10).Value = 4619.2
$ws.Cells.Item(136, 11).Value = 11144.6571
$ws.Cells.Item(136, 12).Value = 13857.6
$ws.Cells.Item(136, 13).Value = -8594.6571
$ws.Cells.Item(136, 14).Value = -18957.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 2787983.8
$ws.Cells.Item(4, 9).Value = 4176610.8
$ws.Cells.Item(4, 10).Value = 10730
$ws.Cells.Item(4, 11).Value = 4176610.8
$ws.Cells.Item(4, 12).Value = 10730
$ws.Cells.Item(4, 13).Value = -4176497.8
$ws.Cells.Item(4, 14).Value = -10956
$ws.Cells.Item(74, 8).Value = 24913.5
$ws.Cells.Item(74, 10).Value = 24913.5
$ws.Cells.Item(74, 12).Value = 24913.5
$ws.Cells.Item(74, 14).Value = -26785.5
$ws.Cells.Item(77, 8).Value = 24913.5
$ws.Cells.Item(77, 10).Value = 24913.5
$ws.Cells.Item(77, 12).Value = 74740.5
$ws.Cells.Item(77, 14).Value = -84100.5
$ws.Cells.Item(132, 8).Value = 2760.48
$ws.Cells.Item(132, 9).Value = 2850.5
$ws.Cells.Item(132, 11).Value = 8551.5
$ws.Cells.Item(132, 13).Value = -6021.5
